# Incorporate typo/wording feedback from coauthors on the "Modeling notes" /
# "Expected relationship with aspen" text in Sheet2.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("C2").Value = "Fast-growing, short-lived species such as aspen typically have high water demand (Ireland et al. 2014). Thus high ADI has been linked with lower climate suitability for aspen (Rehfeldt et al. 2009, 2015)."

$ws2.Range("C3").Value = "To prevent early development of new buds that may lead to injury, aspen phenology is requires a chilling period. Insufficient chilling periods may delay budburst (Man et al. 2017). "

$ws2.Range("C4").Value = "Greater precipitation during the growing season may alleviate summer moisture stress (Worral et al. 2013)."

$ws2.Range("C5").Value = "Plant productivity is higher in areas with warmer temperatures and greater precipitation during the growing season  (Rehfeldt et al. 2009)."

$ws2.Range("C8").Value = "Extreme seasonal variation in temperature may present physiological challenges to aspen (Worrall et al. 2013; Rehfeldt et al. 2015)"

# Reflect the author's last selection before saving.
$ws2.Range("C8").Select()
